# Add a "Date" / temporal-network column to the Journal-articles sheet (Sheet3)
# of the Excel/CSV importer test workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # Sheet3 - the journal-articles table (tabSelected="1")

# --- Column C header + data: "Year" -> "Date" with comma-separated date ranges ---
$ws.Range("C1").Value = "Date"
$ws.Range("C2").Value = "1972-01-01,1972-01-15"
$ws.Range("C3").Value = "1972-01-13,1972-01-19"
$ws.Range("C4").Value = "1972-01-07,1972-01-28"

# Format column C as text (numFmtId 49, i.e. "@") with wrapped text, like the
# other descriptive columns on this sheet.
$dataRange = $ws.Range("C1:C4")
$dataRange.NumberFormat = "@"
$dataRange.WrapText = $true

# Widen column C to fit the new, longer values.
$ws.Columns.Item(3).ColumnWidth = 10.140625

# Rows 2-4 need to grow taller to accommodate the wrapped date-range text.
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 45

# Move the active selection from A5 to C5.
$ws.Range("C5").Select() | Out-Null

# Set up page setup (paper size / orientation) for printing this sheet.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
